# Update cryptos list: refreshed Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.406.37"
$ws.Range("E2").Value = "  +1.92%  "
$ws.Range("D3").Value = "1.884.36"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "0.695"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "'246.70"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "43.52"
$ws.Range("E8").Value = "  +4.87%  "
$ws.Range("E9").Value = "  +3.24%  "
$ws.Range("D10").Value = "0.0746"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").Value = "0.0979"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").Value = "13.55"
$ws.Range("E12").Value = "  +5.42%  "
$ws.Range("D13").Value = "2.160.59"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").Value = "0.772"
$ws.Range("E14").Value = "  +8.35%  "
$ws.Range("D15").Value = "4.95"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "1.870.79"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "35.414.70"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "73.54"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "245.04"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "12.83"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("D22").Value = "5.17"
$ws.Range("E22").Value = "  +5.41%  "
$ws.Range("D23").Value = "2.62"
$ws.Range("E23").Value = "  +9.23%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  -2.65%  "
$ws.Range("D26").Value = "'164.60"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "8.65"
$ws.Range("E27").Value = "  +3.21%  "
$ws.Range("D28").Value = "18.33"
$ws.Range("E28").Value = "  +0.64%  "
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "0.0596"
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").Value = "4.29"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  +4.70%  "
$ws.Range("D33").Value = "4.18"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "1.49"
$ws.Range("E35").Value = "  -9.58%  "
$ws.Range("D36").Value = "0.853"
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("D37").Value = "1.96"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "0.0735"
$ws.Range("E38").Value = "  +11.69%  "
$ws.Range("D39").Value = "17.31"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "0.0218"
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("D41").Value = "'97.40"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "'2.40"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").Value = "1.310.07"
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("D45").Value = "0.0808"
$ws.Range("E45").Value = "  +6.23%  "
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "12.16"
$ws.Range("E48").Value = "  +0.86%  "
$ws.Range("D49").Value = "6.34"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").Value = "42.35"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").Value = "2.065.41"
$ws.Range("E51").Value = "  +0.59%  "

Write-Output "Applied cryptos update"
